$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 19.79629629629628
$ws.Range("R3").Value = 1.911855479578636
$ws.Range("S3").Value = 2.09608909874769
$ws.Range("K16").Value = 21.28240740740739
$ws.Range("R16").Value = 1.940636870984383
$ws.Range("S16").Value = 2.131200751448103
$ws.Range("K17").Value = 21.28240740740739
$ws.Range("K20").Value = 19.30324074074072
$ws.Range("R20").Value = 1.902494195535734
$ws.Range("S20").Value = 2.084694111942012
$ws.Range("K21").Value = 15.74228395061728
$ws.Range("R21").Value = 1.837513876759573
$ws.Range("S21").Value = 2.005936573945218
$ws.Range("K24").Value = 20.22222222222222
$ws.Range("R24").Value = 1.920016703786191
$ws.Range("S24").Value = 2.106033415841584
$ws.Range("K25").Value = 13.75752314814816
$ws.Range("R25").Value = 1.803186500133452
$ws.Range("S25").Value = 1.964569140204562
$ws.Range("K28").Value = 12.93898809523811
$ws.Range("R28").Value = 1.789400236291612
$ws.Range("S28").Value = 1.948001533154466
$ws.Range("K29").Value = 12.93898809523811
$ws.Range("K34").Value = 5.486111111111112
$ws.Range("R34").Value = 1.672941176470588
$ws.Range("S34").Value = 1.809089700996678
$ws.Range("K35").Value = 21.28240740740739
$ws.Range("R35").Value = 1.940636870984383
$ws.Range("S35").Value = 2.131200751448103
$ws.Range("K39").Value = 1.791666666666668
$ws.Range("R39").Value = 1.620655622136059
$ws.Range("S39").Value = 1.747323835194455
$ws.Range("K40").Value = 5.486111111111112
$ws.Range("R40").Value = 1.672941176470588
$ws.Range("S40").Value = 1.809089700996678
$ws.Range("K44").Value = 19.30324074074072
$ws.Range("R44").Value = 1.902494195535734
$ws.Range("S44").Value = 2.084694111942012
$ws.Range("K46").Value = 12.67039049919483
$ws.Range("R46").Value = 1.784922174701128
$ws.Range("S46").Value = 1.942625691911729
$ws.Range("K54").Value = 21.28240740740739
$ws.Range("R54").Value = 1.940636870984383
$ws.Range("S54").Value = 2.131200751448103
$ws.Range("K57").Value = 13.62268518518517
$ws.Range("R57").Value = 1.80090088129692
$ws.Range("S57").Value = 1.961820583643568
$ws.Range("K60").Value = 12.67039049919483
$ws.Range("R60").Value = 1.784922174701128
$ws.Range("S60").Value = 1.942625691911729
$ws.Range("K62").Value = 12.67039049919483
$ws.Range("R62").Value = 1.784922174701128
$ws.Range("S62").Value = 1.942625691911729
$ws.Range("K65").Value = 19.65277777777778
$ws.Range("R65").Value = 1.909121107266436
$ws.Range("S65").Value = 2.092759415833974
$ws.Range("K73").Value = 19.30324074074072
$ws.Range("R73").Value = 1.902494195535734
$ws.Range("S73").Value = 2.084694111942012
$ws.Range("K86").Value = 13.76976495726495
$ws.Range("R86").Value = 1.803394296576035
$ws.Range("S86").Value = 1.964819060413116
$ws.Range("K89").Value = 21.28240740740739
$ws.Range("R89").Value = 1.940636870984383
$ws.Range("S89").Value = 2.131200751448103
$ws.Range("K90").Value = 13.75752314814816
$ws.Range("K91").Value = 13.75752314814816
$ws.Range("R91").Value = 1.803186500133452
$ws.Range("S91").Value = 1.964569140204562
